# agregué columna nueva usdm2
# Adds three new documentation rows (43-45) describing the new columns
# usd_per_m2_2, provincia_ciudad and ciudad_barrio, mirroring the existing
# "field name / dtype / description" triples already present in columns D/F/I.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 43: usd_per_m2_2 / float64 / 17817 no nulos
$ws.Range("D43").Value2 = "usd_per_m2_2"
$ws.Range("F43").Value2 = "float64"
$ws.Range("I43").Value2 = "17817 no nulos"

# Row 44: provincia_ciudad / str / description
$ws.Range("D44").Value2 = "provincia_ciudad"
$ws.Range("F44").Value2 = "str"
$ws.Range("I44").Value2 = "provincia si es el interior - capital - bsas zona sur o zona norte"

# Row 45: ciudad_barrio / str / description
$ws.Range("D45").Value2 = "ciudad_barrio"
$ws.Range("F45").Value2 = "str"
$ws.Range("I45").Value2 = "ciudad si es el interior - barrio si es capital"

# Move the selection / viewport the way the author left the sheet: scrolled
# down so row 10 is at the top, with I46 as the active cell.
$ws.Range("I46").Select()
$excel.ActiveWindow.ScrollRow = 10
$excel.ActiveWindow.ScrollColumn = 1
